$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 339, shifting the existing rows 339-360 down to 340-361
$ws.Rows.Item(339).Insert()

# Populate the newly inserted row 339 with the new weekly record
$ws.Cells.Item(339, 1).Value = 7
$ws.Cells.Item(339, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(339, 3).Value = "Ñuble"
$ws.Cells.Item(339, 4).Value = 45265
$ws.Cells.Item(339, 5).Value = 16
$ws.Cells.Item(339, 6).Value = 100112040
$ws.Cells.Item(339, 7).Value = "Cilantro"
$ws.Cells.Item(339, 8).Value = "Sin especificar"
$ws.Cells.Item(339, 9).Value = "Primera"
$ws.Cells.Item(339, 10).Value = 240
$ws.Cells.Item(339, 11).Value = 1500
$ws.Cells.Item(339, 12).Value = 2000
$ws.Cells.Item(339, 13).Value = 1750
$ws.Cells.Item(339, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(339, 15).Value = "Región de Ñuble"
$ws.Cells.Item(339, 16).Value = 1750
$ws.Cells.Item(339, 17).Value = 1
$ws.Cells.Item(339, 18).Value = "Hortaliza"

# Make sure the date cell keeps the workbook's date number format
$ws.Cells.Item(339, 4).NumberFormat = $ws.Cells.Item(340, 4).NumberFormat
